# Actualizacion automatica del mapa (2026-02-24 14:38:08)
#
# Source data updates the PEBCOM case list:
#   1. Removes the now-closed case "-747" (Ugarteche 2816) that used to sit
#      at row 201 -- every following row shifts up by one.
#   2. Appends 8 freshly-reported cases at the bottom of the table
#      (rows 214-221 once the sheet has settled after the removal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the obsolete case row; Excel shifts everything below it up ---
$ws.Rows.Item(201).Delete()

# --- 2. New cases to append at the end of the table ---
$newCases = @(
    @{ Caso="-749";      Fecha="2/23/2026"; Direccion="Moreno 2965";                           Comuna="3";  Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="base corroida ";              Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.407758; Y=-34.613793; Operacion="Almagro";   Zona="Capital Sur"; PD="CEN-I"; N2="Fuera de Poligono OVL" }
    @{ Caso="-752";      Fecha="2/23/2026"; Direccion="RIVADAVIA AV. 5691";                     Comuna="6";  Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="base corroida";               Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.445685; Y=-34.622144; Operacion="Boedo";     Zona="Capital Sur"; PD="PCH-G"; N2="Fuera de Poligono OVL" }
    @{ Caso="-756";      Fecha="2/23/2026"; Direccion="SINCLAIR 3106";                          Comuna="14"; Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="inclinada";                   Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.422892; Y=-34.573802; Operacion="Palermo";   Zona="Capital Sur"; PD="VCR-M"; N2="Fuera de Poligono OVL" }
    @{ Caso="-760";      Fecha="2/23/2026"; Direccion="MAZA 1615";                              Comuna="5";  Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="base corroida e inclinada";   Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.414586; Y=-34.630095; Operacion="Boedo";     Zona="Capital Sur"; PD="PPT-E"; N2="ARATO-25058.PO.1PPT" }
    @{ Caso="S00519068"; Fecha="2/23/2026"; Direccion="CALVO, CARLOS 3747";                     Comuna="5";  Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="base corroida";               Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.41859;  Y=-34.624508; Operacion="Boedo";     Zona="Capital Sur"; PD="ALM-A"; N2="Fuera de Poligono OVL" }
    @{ Caso="S01061920"; Fecha="2/23/2026"; Direccion="BROWN, ALTE. AV. 1184 ";                 Comuna="4";  Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="corroida inclinada";          Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.358913; Y=-34.635093; Operacion="San Telmo"; Zona="Capital Sur"; PD="CON-G"; N2="Fuera de Poligono OVL" }
    @{ Caso="S01064368"; Fecha="2/23/2026"; Direccion=" SCALABRINI ORTIZ, RAUL AV. 1413";       Comuna="14"; Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="inclinada";                   Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.426552; Y=-34.592076; Operacion="Palermo";   Zona="Capital Sur"; PD="VCR-D"; N2="Fuera de Poligono OVL" }
    @{ Caso="-764";      Fecha="2/24/2026"; Direccion="CASEROS AV. 3547";                       Comuna="5";  Estado="Pendiente"; Proveedor="PEBCOM"; OT="Pendiente ADM"; Obs="colocar en 3543";             Tarea="Cambio"; Equipo="Sin equipos"; Elemento="Pasante"; Att=1; X=-58.41515;  Y=-34.638705; Operacion="San Telmo"; Zona="Capital Sur"; PD="PPT-Q"; N2="Fuera de Poligono OVL" }
)

# First free row right after the current used range
$startRow = $ws.UsedRange.Rows.Count + 1

$r = $startRow
foreach ($case in $newCases) {

    # Columns A-K and O-R hold free-form/identifier text in this sheet (case
    # numbers, comuna codes, dates, etc. are all stored as text, never as
    # numbers/dates) -- force text formatting before writing so lookalike
    # values such as "-749", "3" or "2/23/2026" aren't reinterpreted by Excel.
    $textRange = $ws.Range("A" + $r + ":K" + $r + ",O" + $r + ":R" + $r)
    $textRange.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value  = $case.Caso
    $ws.Cells.Item($r, 2).Value  = $case.Fecha
    $ws.Cells.Item($r, 3).Value  = $case.Direccion
    $ws.Cells.Item($r, 4).Value  = $case.Comuna
    $ws.Cells.Item($r, 5).Value  = $case.Estado
    $ws.Cells.Item($r, 6).Value  = $case.Proveedor
    $ws.Cells.Item($r, 7).Value  = $case.OT
    $ws.Cells.Item($r, 8).Value  = $case.Obs
    $ws.Cells.Item($r, 9).Value  = $case.Tarea
    $ws.Cells.Item($r, 10).Value = $case.Equipo
    $ws.Cells.Item($r, 11).Value = $case.Elemento
    $ws.Cells.Item($r, 15).Value = $case.Operacion
    $ws.Cells.Item($r, 16).Value = $case.Zona
    $ws.Cells.Item($r, 17).Value = $case.PD
    $ws.Cells.Item($r, 18).Value = $case.N2

    # Back to the sheet's default (unstyled) look now that the text is safe
    $textRange.Style = "Normal"

    # Attachments / coordinates are genuine numbers in every other row
    $ws.Cells.Item($r, 12).Value = $case.Att
    $ws.Cells.Item($r, 13).Value = $case.X
    $ws.Cells.Item($r, 14).Value = $case.Y

    $r++
}
